$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.336.56'
$ws.Range('E2').Value = '  +3.45%  '
$ws.Range('D3').Value = '3.219.69'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '540.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.61'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.89%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.539'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.50%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('E10').Value = '  +4.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.438'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.70%  '
$ws.Range('D12').Value = '3.773.42'
$ws.Range('E12').Value = '  +2.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.138'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('E15').Value = '  +3.13%  '
$ws.Range('D16').Value = '60.378.95'
$ws.Range('E16').Value = '  +3.44%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.179.29'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '382.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('E23').Value = '  +2.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.170'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').Value = '0.0₃0912'
$ws.Range('E28').Value = '  +3.65%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.91'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.51'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.46'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.50%  '
$ws.Range('E33').Value = '  +4.05%  '
$ws.Range('E34').Value = '  +6.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '157.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.51%  '
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('B37').Value = 'EnergySwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.03'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '2.797.28'
$ws.Range('E38').Value = '  +5.58%  '
$ws.Range('E39').Value = '  +4.51%  '
$ws.Range('E40').Value = '  +0.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.28'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '40.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.96%  '
$ws.Range('E43').Value = '  +2.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0290'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.73%  '
$ws.Range('D45').Value = '3.261.15'
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('E46').Value = '  +3.13%  '
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.808'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.34%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.87'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '274.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.83%  '
